# Generate Report for Handoff
# Update the status/handoff-date for the 82b901f4-3ba9-4ba9-b5ca-a91d40e19d36
# file now that it is ready for handoff (rather than "Handed back: in sync
# with en-US").

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 10:52:43"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-23 10:52:39"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-23 10:52:43"
